$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - Research Development
$ws.Range("D17").Value = 2444700

# Row 18 - Selling General and Administrative
$ws.Range("D18").Value = 475700

# Row 21 - Total Operating Expenses
$ws.Range("D21").Value = 599700
$ws.Range("J21").Value = "NA"

# Row 23 - Earnings Before Interest And Taxes
$ws.Range("D23").Value = 478000

# Row 24 - Interest Expense
$ws.Range("D24").Value = 93100

# Row 26 - Income After Tax
$ws.Range("D26").Value = 384900

# Row 27 - Net Income From Continuing Ops
$ws.Range("D27").Value = 384900

# Row 33 - Net Income
$ws.Range("D33").Value = 261900

# Row 35 - Net Income Applicable To Common Shares
$ws.Range("D35").Value = 261900

# Row 58 - Other Current Liabilities
$ws.Range("D58").Value = 23800
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "NA"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"

# Row 59 - Total Current Liabilities
$ws.Range("D59").Value = 2225100
$ws.Range("E59").Value = 1814900
$ws.Range("F59").Value = 1495900
$ws.Range("G59").Value = 1303000
$ws.Range("H59").Value = 1129800
$ws.Range("I59").Value = 959100

# Row 61 - Other Liabilities
$ws.Range("E61").Value = 745600

# Row 62 - Total Liabilities
$ws.Range("D62").Value = 1151900
$ws.Range("E62").Value = 651200

# Row 81 - Capital Expenditures
$ws.Range("D81").Value = 261900

# Row 83 - Depreciation
$ws.Range("J83").Value = "NA"

# Row 91 - Changes In Other Operating Activities
$ws.Range("D91").Value = -85000
$ws.Range("E91").Value = -69100
$ws.Range("F91").Value = -41600
$ws.Range("G91").Value = -45600
$ws.Range("H91").Value = -79600
$ws.Range("I91").Value = -85700
$ws.Range("J91").Value = -46300

# Row 94 - Total Cash Flows From Investing Activities
$ws.Range("J94").Value = "NA"

# Row 100 - Total Cash Flows From Financing Activities
$ws.Range("J100").Value = "NA"

# Row 101 - Effect Of Exchange Rate Changes
$ws.Range("J101").Value = "NA"
